# Automatic update of files.
# Rows 13-15 are cyclically rotated (13<-14, 14<-15, 15<-13) and rows 22-23 are swapped,
# mirroring the observation re-ordering captured in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually differ between the affected rows (rest of each row -
# C, I, K, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY - is identical
# across the group and is left untouched).
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")

function Get-RowData($row) {
    $data = @{}
    foreach ($col in $cols) {
        $data[$col] = $ws.Range("$col$row").Value2
    }
    $data["M"] = $ws.Range("M$row").Value2
    return $data
}

function Set-RowData($row, $data) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $data[$col]
    }
    if ($data["M"]) {
        $ws.Range("M$row").Value = $data["M"]
    } else {
        $ws.Range("M$row").Value = ""
    }
}

# --- Rows 13, 14, 15: cyclic rotation up (new13 = old14, new14 = old15, new15 = old13)
$old13 = Get-RowData 13
$old14 = Get-RowData 14
$old15 = Get-RowData 15

Set-RowData 13 $old14
Set-RowData 14 $old15
Set-RowData 15 $old13

# --- Rows 22, 23: swap
$old22 = Get-RowData 22
$old23 = Get-RowData 23

Set-RowData 22 $old23
Set-RowData 23 $old22
